$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.085.72"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.60"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.30"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.90"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("E7").Value = "  +3.13%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.77"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.86"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.95"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.673.13"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.442.86"
$ws.Range("E16").Value = "  +5.91%  "
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.997.46"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.51"
$ws.Range("E19").Value = "  +6.80%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.20"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.89"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.93"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.57"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  -9.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.18"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.36"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.93"
$ws.Range("E32").Value = "  +6.07%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.22"
$ws.Range("E33").Value = "  +3.69%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.39"
$ws.Range("E35").Value = "  +7.76%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0696"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.997.32"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0288"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.47"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.83"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.75"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.539.52"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.07"
$ws.Range("E50").Value = "  +5.28%  "
$ws.Range("E51").Value = "  +0.85%  "
